$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 467, pushing existing rows 467:496 down to 468:497.
$ws.Rows(467).Insert()

# Populate the newly inserted row 467 with the new record.
# (Columns A,B,C,E,F,G,H,I,J,K,L,M,N,P,Q,R mirror the record that used to be
# at the bottom of the sheet; D (Fecha) and O (Origen) are new values.)
$ws.Cells.Item(467, 1).Value = 9
$ws.Cells.Item(467, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(467, 3).Value = "Metropolitana"
$ws.Cells.Item(467, 4).Value = 44931
$ws.Cells.Item(467, 5).Value = 13
$ws.Cells.Item(467, 6).Value = 100112032
$ws.Cells.Item(467, 7).Value = "Zapallo italiano"
$ws.Cells.Item(467, 8).Value = "Sin especificar"
$ws.Cells.Item(467, 9).Value = "Primera"
$ws.Cells.Item(467, 10).Value = 340
$ws.Cells.Item(467, 11).Value = 4000
$ws.Cells.Item(467, 12).Value = 5000
$ws.Cells.Item(467, 13).Value = 4500
$ws.Cells.Item(467, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(467, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(467, 16).Value = 90
$ws.Cells.Item(467, 17).Value = 50
$ws.Cells.Item(467, 18).Value = "Hortaliza"

# Make sure the Fecha (date) column keeps the date number format used by the
# rest of the column.
$ws.Cells.Item(467, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
